$p = $ppt.ActivePresentation

# --- Slide 1: "Rectangle 3" subtitle placeholder ---
# Merge the runs "IIR.C11" and ". Probabilistic " into a single run
# "IIR.C11. Probabilistic " (leaving "information retrieval" untouched).
$s1 = $p.Slides.Item(1)
$shp1 = $s1.Shapes.Item(2)
$tr1 = $shp1.TextFrame.TextRange
$para2 = $tr1.Paragraphs(2, 1)
$mergeRange = $para2.Characters(1, 23)
$mergeRange.Text = "IIR.C11. Probabilistic "

# --- Slide 36: title placeholder ---
# Change "Bài tập" to "Bài tập 5.1" (inserting " 5.1" right after "tập").
$s36 = $p.Slides.Item(36)
$shp36 = $s36.Shapes.Item(1)
$tr36 = $shp36.TextFrame.TextRange
$editRange = $tr36.Characters(5, 3)
$editRange.Text = "tập 5.1"
